# Update countries & provincias Spain
#
# Refreshes COVID-19 case figures for several countries and keeps the
# "Pais" sheet sorted by "Casos totales" (column B) descending, which
# means a few rows swap places with their neighbours once the new
# totals are in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Simple data refreshes (country keeps its current rank) -----------

Set-Row 31 @("Israel", 16237, 29, 9858, 6145, 103, 2, 234)
Set-Row 53 @("Malasia", 6353, 55, 4484, 1764, 28, 0, 105)
Set-Row 92 @("Hong Kong", 1041, 1, 900, 137, 1, 0, 4)

# --- Banglades overtakes Dinamarca / Filipinas / Serbia (rows 41-44) --
# Banglades' new total (10143) pushes it above Dinamarca, Filipinas and
# Serbia, which each shift down one row, keeping their prior figures.

Set-Row 41 @("Banglades", 10143, 688, 1209, 8752, 1, 5, 182)
Set-Row 42 @("Dinamarca", 9670, 147, 6987, 2199, 62, 0, 484)
Set-Row 43 @("Filipinas", 9485, 262, 1315, 7547, 31, 16, 623)
Set-Row 44 @("Serbia", 9464, 0, 1551, 7720, 54, 0, 193)

# --- Benin overtakes Monaco (rows 156-157) -----------------------------

Set-Row 156 @("Benin", 96, 6, 50, 44, 0, 0, 2)
Set-Row 157 @("Monaco", 95, 0, 78, 13, 1, 0, 4)

# --- Belice / Santa Lucia swap places (rows 188-189) -------------------

Set-Row 188 @("Belice", 18, 0, 13, 3, 1, 0, 2)
Set-Row 189 @("Santa Lucia", 18, 0, 15, 3, 0, 0, 0)

# --- San Vicente y las Granadinas / Namibia swap (rows 194-195) --------

Set-Row 194 @("San Vicente y las Granadinas", 16, 0, 8, 8, 0, 0, 0)
Set-Row 195 @("Namibia", 16, 0, 8, 8, 0, 0, 0)

# --- San Cristobal y Nieves / Burundi swap (rows 198-199) --------------

Set-Row 198 @("San Cristobal y Nieves", 15, 0, 8, 7, 0, 0, 0)
Set-Row 199 @("Burundi", 15, 0, 7, 7, 0, 0, 1)
